$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("games")

# Update the bp1-titka game row (row 2): age, type, title.
# Write in this order so the shared-string table grows in the same
# sequence as the target workbook (age "14+" first, then type "normal",
# then the new title).
$ws.Range("D2").Value = "14+"
$ws.Range("B2").Value = "normal"
$ws.Range("C2").Value = "BudapestOne az üvegfalak titka"

# Make the "games" sheet the active/selected tab, with C2 selected,
# scrolled back to the top-left (A1), replacing "tasks" as active sheet.
$ws.Activate()
$ws.Range("C2").Select()
